# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") values recomputed; write the new literal values row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 7
    3  = 5
    4  = 11
    5  = 5
    6  = 1
    7  = 10
    8  = 7
    9  = 3
    10 = 8
    11 = 3
    12 = 2
    13 = 6
    14 = 5
    15 = 0
    16 = 7
    17 = 2
    18 = 9
    19 = 8
    20 = 3
    21 = 3
    22 = 6
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
